$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top (shifts existing data down)
$ws.Rows.Item(1).Insert()

# Add header labels
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Value"

# Update the selected cell to C1, matching the diff
$ws.Range("C1").Select()
